# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# header style of the existing last header cell (H1), and fill in the
# per-row values for rows 2..71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style used by the other header cells (bold, centered, bordered)
# by copying the formatting from the existing H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:I71 and J2:J71
$i0Values = @(9,7,5,6,7,6,9,7,9,8,8,7,7,7,7,7,6,8,8,5,4,7,6,6,6,6,8,4,6,7,3,7,9,7,5,7,10,4,6,6,3,5,8,8,5,8,6,1,7,11,4,9,6,7,1,7,6,6,5,5,6,8,9,8,6,5,5,4,9,9)
$ifValues = @(9,8,5,6,8,6,9,7,9,8,8,7,8,8,7,7,6,9,9,6,5,7,6,6,6,6,8,5,6,7,5,7,9,7,5,8,10,6,6,7,4,7,8,8,6,8,6,2,8,11,5,10,6,8,3,7,6,6,5,7,6,8,10,8,7,5,5,4,9,9)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
